$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition) - sheetId 1
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 7299
$ws1.Range("F3").Value = 63
$ws1.Range("F5").Value = 179
$ws1.Range("F6").Value = 1109
$ws1.Range("F7").Value = 188
$ws1.Range("F8").Value = 10
$ws1.Range("F9").Value = 92
$ws1.Range("F10").Value = 25

# Sheet "演出" (Performance) - sheetId 2
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 12

# Sheet "全部类型" (All types) - sheetId 4
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 7299
$ws4.Range("F3").Value = 63
$ws4.Range("F5").Value = 179
$ws4.Range("F6").Value = 1109
$ws4.Range("F7").Value = 188
$ws4.Range("F8").Value = 12
$ws4.Range("F9").Value = 10
$ws4.Range("F10").Value = 92
$ws4.Range("F11").Value = 25
